$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 39, shifting existing rows 39-45 down to 40-46
$ws.Rows.Item(39).Insert()

# Match the style used by the date column (D) in the other rows
$ws.Cells.Item(39, 4).NumberFormat = $ws.Cells.Item(40, 4).NumberFormat

# Populate new row 39 with data
$ws.Cells.Item(39, 1).Value = 11
$ws.Cells.Item(39, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(39, 3).Value = "Bíobío"
$ws.Cells.Item(39, 4).Value = 44748
$ws.Cells.Item(39, 5).Value = 8
$ws.Cells.Item(39, 6).Value = 100112013
$ws.Cells.Item(39, 7).Value = "Alcachofa"
$ws.Cells.Item(39, 8).Value = "Argentina(o)"
$ws.Cells.Item(39, 9).Value = "Primera"
$ws.Cells.Item(39, 10).Value = 110
$ws.Cells.Item(39, 11).Value = 15000
$ws.Cells.Item(39, 12).Value = 16000
$ws.Cells.Item(39, 13).Value = 15545
$ws.Cells.Item(39, 14).Value = "$/caja 40 unidades"
$ws.Cells.Item(39, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(39, 16).Value = 389
$ws.Cells.Item(39, 17).Value = 40
$ws.Cells.Item(39, 18).Value = "Hortaliza"
